$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.335.76"
$ws.Range("E2").Value = "  +1.09%  "
$ws.Range("D3").Value = "1.621.96"
$ws.Range("E3").Value = "  +1.52%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Value = "'212.02"
$ws.Range("E5").Value = "  +0.03%  "
$ws.Range("E6").Value = "  +0.08%  "
$ws.Range("E7").Value = "  +0.40%  "
$ws.Range("E8").Value = "  -0.22%  "
$ws.Range("E9").Value = "  +0.22%  "
$ws.Range("D10").Value = "'18.67"
$ws.Range("E10").Value = "  +2.04%  "
$ws.Range("E11").Value = "  +0.65%  "
$ws.Range("D12").Value = "1.848.41"
$ws.Range("E12").Value = "  +1.58%  "
$ws.Range("D13").Value = "1.626.30"
$ws.Range("E13").Value = "  +1.73%  "
$ws.Range("E14").Value = "  +0.24%  "
$ws.Range("E15").Value = "  +0.08%  "
$ws.Range("D16").Value = "26.347.85"
$ws.Range("E16").Value = "  +1.13%  "
$ws.Range("D17").Value = "'62.39"
$ws.Range("E17").Value = "  +2.48%  "
$ws.Range("D18").Value = "0.0₃0724"
$ws.Range("E18").Value = "  -0.57%  "
$ws.Range("E19").Value = "  +0.04%  "
$ws.Range("D20").Value = "'202.26"
$ws.Range("E20").Value = "  -0.78%  "
$ws.Range("E21").Value = "  -0.16%  "
$ws.Range("D22").Value = "'9.30"
$ws.Range("E22").Value = "  +0.26%  "
$ws.Range("D23").Value = "'6.04"
$ws.Range("E23").Value = "  +0.23%  "
$ws.Range("E24").Value = "  -2.29%  "
$ws.Range("D25").Value = "'144.63"
$ws.Range("E25").Value = "  +0.60%  "
$ws.Range("E26").Value = "  +0.10%  "
$ws.Range("D27").Value = "'0.119"
$ws.Range("E27").Value = "  -1.25%  "
$ws.Range("D28").Value = "'15.19"
$ws.Range("E28").Value = "  -0.28%  "
$ws.Range("E29").Value = "  +0.47%  "
$ws.Range("D30").Value = "'0.0518"
$ws.Range("E30").Value = "  +8.46%  "
$ws.Range("E31").Value = "  +0.02%  "
$ws.Range("E32").Value = "  +1.16%  "
$ws.Range("E33").Value = "  +0.43%  "
$ws.Range("E34").Value = "  -0.21%  "
$ws.Range("D36").Value = "1.156.93"
$ws.Range("E36").Value = "  +2.03%  "
$ws.Range("E37").Value = "  +0.47%  "
$ws.Range("D38").Value = "'0.804"
$ws.Range("E38").Value = "  +0.67%  "
$ws.Range("E39").Value = "  +0.11%  "
$ws.Range("E40").Value = "  +0.08%  "
$ws.Range("D41").Value = "'0.497"
$ws.Range("E41").Value = "  +0.47%  "
$ws.Range("E42").Value = "  +3.83%  "
$ws.Range("D43").Value = "'0.783"
$ws.Range("E43").Value = "  +0.55%  "
$ws.Range("D44").Value = "1.759.48"
$ws.Range("E44").Value = "  +1.50%  "
$ws.Range("D45").Value = "'92.60"
$ws.Range("E45").Value = "  +0.56%  "
$ws.Range("B46").Value = "RenderToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D46").Value = "'1.52"
$ws.Range("E46").Value = "  +0.83%  "
$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").Value = "'53.76"
$ws.Range("E47").Value = "  -0.78%  "
$ws.Range("B48").Value = "Cronos"
$ws.Range("C48").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D48").Value = "'0.0509"
$ws.Range("E48").Value = "  +0.63%  "
$ws.Range("B49").Value = "Mantle"
$ws.Range("C49").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D49").Value = "'0.409"
$ws.Range("E49").Value = "  +0.63%  "
$ws.Range("B50").Value = "USDD"
$ws.Range("C50").Value = "https://coinranking.com/coin/z2PZIKQL7+usdd-usdd"
$ws.Range("D50").Value = "'1.00"
$ws.Range("E50").Value = "  -0.24%  "
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").Value = "'7.30"
$ws.Range("E51").Value = "  +2.01%  "
